# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 01:35"

# Row 4 - Estados Unidos (idx 8) - totals updated
$ws.Range("B4").Value = 1619798
$ws.Range("C4").Value = 27075
$ws.Range("D4").Value = 381923
$ws.Range("E4").Value = 1141605
$ws.Range("G4").Value = 1334
$ws.Range("H4").Value = 96270

# Rows 61/62 - Nigeria overtakes Moldavia in case count, so they swap rows.
# Row 61 becomes Nigeria with updated figures.
$ws.Range("A61").Value = "Nigeria"
$ws.Range("B61").Value = 7016
$ws.Range("C61").Value = 339
$ws.Range("D61").Value = 1907
$ws.Range("E61").Value = 4898
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 211

# Row 62 becomes Moldavia with its previous (unchanged) figures.
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 6704
$ws.Range("C62").Value = 151
$ws.Range("D62").Value = 2953
$ws.Range("E62").Value = 3518
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 233

# Row 103
$ws.Range("D103").Value = 375
$ws.Range("E103").Value = 684

# Row 120
$ws.Range("B120").Value = 749
$ws.Range("C120").Value = 3
$ws.Range("D120").Value = 594
$ws.Range("E120").Value = 135

# Row 180
$ws.Range("B180").Value = 51
$ws.Range("C180").Value = 3
$ws.Range("E180").Value = 29
